$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2391.1
$ws.Range("I28").Value = 1484.3334
$ws.Range("J28").Value = 3751.25
$ws.Range("K28").Value = 1484.3334
$ws.Range("L28").Value = 3751.25
$ws.Range("M28").Value = -999.3334
$ws.Range("N28").Value = -4721.25

$ws.Range("H41").Value = 481.3846
$ws.Range("J41").Value = 460.33334
$ws.Range("L41").Value = 460.33334
$ws.Range("N41").Value = -1340.33334

$ws.Range("H53").Value = 3125.1765
$ws.Range("I53").Value = 1394.3334
$ws.Range("K53").Value = 1394.3334
$ws.Range("M53").Value = -757.3334

$ws.Range("H76").Value = 6031.4614
$ws.Range("I76").Value = 5260.8184
$ws.Range("J76").Value = 6596.6
$ws.Range("K76").Value = 5260.8184
$ws.Range("L76").Value = 6596.6
$ws.Range("M76").Value = -4945.8184
$ws.Range("N76").Value = -7226.6

$ws.Range("H79").Value = 6031.4614
$ws.Range("I79").Value = 5260.8184
$ws.Range("J79").Value = 6596.6
$ws.Range("K79").Value = 5260.8184
$ws.Range("L79").Value = 6596.6
$ws.Range("M79").Value = -4168.8184
$ws.Range("N79").Value = -8780.6

$ws.Range("H86").Value = 5692.5386
$ws.Range("I86").Value = 7501.5
$ws.Range("J86").Value = 5363.636
$ws.Range("K86").Value = 7501.5
$ws.Range("L86").Value = 5363.636
$ws.Range("M86").Value = -6378.5
$ws.Range("N86").Value = -7609.636

$ws.Range("H89").Value = 5692.5386
$ws.Range("I89").Value = 7501.5
$ws.Range("J89").Value = 5363.636
$ws.Range("K89").Value = 37507.5
$ws.Range("L89").Value = 26818.18
$ws.Range("M89").Value = -31891.5
$ws.Range("N89").Value = -38050.18

$ws.Range("H106").Value = 2965.4644
$ws.Range("I106").Value = 1751.9445
$ws.Range("J106").Value = 5149.8
$ws.Range("K106").Value = 1751.9445
$ws.Range("L106").Value = 5149.8
$ws.Range("M106").Value = -1120.9445
$ws.Range("N106").Value = -6411.8

$ws.Range("H107").Value = 338
$ws.Range("I107").Value = 357.33334
$ws.Range("J107").Value = 280
$ws.Range("K107").Value = 357.33334
$ws.Range("L107").Value = 280
$ws.Range("M107").Value = 1562.66666
$ws.Range("N107").Value = -4120

$ws.Range("H112").Value = 1783
$ws.Range("J112").Value = 1834.7142
$ws.Range("L112").Value = 5504.142599999999
$ws.Range("N112").Value = -7720.142599999999

$ws.Range("H132").Value = 798.1316
$ws.Range("I132").Value = 680.82855
$ws.Range("K132").Value = 2042.48565
$ws.Range("M132").Value = 487.5143500000001

$ws.Range("H138").Value = 2996.3132
$ws.Range("I138").Value = 1994.375
$ws.Range("J138").Value = 3235.582
$ws.Range("K138").Value = 5983.125
$ws.Range("L138").Value = 9706.745999999999
$ws.Range("M138").Value = -843.125
$ws.Range("N138").Value = -19986.746

$ws.Range("H141").Value = 2246.625
$ws.Range("I141").Value = 2246.625
$ws.Range("K141").Value = 6739.875
$ws.Range("M141").Value = -1559.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5146.25
$ws.Range("I32").Value = 3940.2222
$ws.Range("J32").Value = 16000.5
$ws.Range("K32").Value = 3940.2222
$ws.Range("L32").Value = 16000.5
$ws.Range("M32").Value = -3653.2222
$ws.Range("N32").Value = -16574.5

$ws.Range("H45").Value = 47620496
$ws.Range("I45").Value = 62501420
$ws.Range("J45").Value = 1541
$ws.Range("K45").Value = 62501420
$ws.Range("L45").Value = 1541
$ws.Range("M45").Value = -62501043
$ws.Range("N45").Value = -2295

$ws.Range("H102").Value = 2750.25
$ws.Range("I102").Value = 2750.25
$ws.Range("K102").Value = 2750.25
$ws.Range("M102").Value = -1128.25

$ws.Range("H110").Value = 5691.357
$ws.Range("I110").Value = 5024.864
$ws.Range("J110").Value = 8135.1665
$ws.Range("K110").Value = 5024.864
$ws.Range("L110").Value = 8135.1665
$ws.Range("M110").Value = -2979.864
$ws.Range("N110").Value = -12225.1665

$ws.Range("H122").Value = 3335.4055
$ws.Range("I122").Value = 2672.2593
$ws.Range("K122").Value = 8016.777900000001
$ws.Range("M122").Value = -5566.777900000001

$ws.Range("H132").Value = 6110.073
$ws.Range("I132").Value = 6353.5
$ws.Range("J132").Value = 5446.1816
$ws.Range("K132").Value = 19060.5
$ws.Range("L132").Value = 16338.5448
$ws.Range("M132").Value = -16530.5
$ws.Range("N132").Value = -21398.5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""

$ws.Range("H134").Value = 3540.4
$ws.Range("I134").Value = 2048.75
$ws.Range("K134").Value = 6146.25
$ws.Range("M134").Value = -3611.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27592.238
$ws.Range("I31").Value = 1899.9166
$ws.Range("K31").Value = 1899.9166
$ws.Range("M31").Value = -1604.9166

$ws.Range("H34").Value = 27592.238
$ws.Range("I34").Value = 1899.9166
$ws.Range("K34").Value = 1899.9166
$ws.Range("M34").Value = -1697.9166

$ws.Range("H105").Value = 3336.6667
$ws.Range("I105").Value = 1158.3334
$ws.Range("K105").Value = 1158.3334
$ws.Range("M105").Value = 588.6666

$ws.Range("H134").Value = 3699.6843
$ws.Range("I134").Value = 2366.6924
$ws.Range("J134").Value = 6587.8335
$ws.Range("K134").Value = 7100.0772
$ws.Range("L134").Value = 19763.5005
$ws.Range("M134").Value = -4565.0772
$ws.Range("N134").Value = -24833.5005

$ws.Range("H135").Value = 69597.75
$ws.Range("J135").Value = 69597.75
$ws.Range("L135").Value = 69597.75
$ws.Range("N135").Value = -79737.75

$ws.Range("H141").Value = 528952.25
$ws.Range("J141").Value = 688603
$ws.Range("L141").Value = 688603
$ws.Range("N141").Value = -698963

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 345000.6
$ws.Range("J37").Value = 345000.6
$ws.Range("L37").Value = 1035001.8
$ws.Range("N37").Value = -1035225.8

$ws.Range("H74").Value = 13009.667
$ws.Range("I74").Value = 999
$ws.Range("J74").Value = 19015
$ws.Range("K74").Value = 2997
$ws.Range("L74").Value = 57045
$ws.Range("M74").Value = -1936
$ws.Range("N74").Value = -59167

$ws.Range("H77").Value = 13009.667
$ws.Range("I77").Value = 999
$ws.Range("J77").Value = 19015
$ws.Range("K77").Value = 8991
$ws.Range("L77").Value = 171135
$ws.Range("M77").Value = -3687
$ws.Range("N77").Value = -181743

$ws.Range("H95").Value = 15027
$ws.Range("J95").Value = 15027
$ws.Range("L95").Value = 45081
$ws.Range("N95").Value = -49199

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 29990
$ws.Range("J32").Value = 29990
$ws.Range("L32").Value = 29990
$ws.Range("N32").Value = -30582

$ws.Range("H42").Value = 37994
$ws.Range("J42").Value = 37994
$ws.Range("L42").Value = 37994
$ws.Range("N42").Value = -38964

$ws.Range("H115").Value = 37994
$ws.Range("J115").Value = 37994
$ws.Range("L115").Value = 37994
$ws.Range("N115").Value = -40344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10217.363
$ws.Range("I7").Value = 6848.6665
$ws.Range("J7").Value = 14259.8
$ws.Range("K7").Value = 6848.6665
$ws.Range("L7").Value = 14259.8
$ws.Range("M7").Value = -6736.6665
$ws.Range("N7").Value = -14483.8

$ws.Range("H22").Value = 2394.742
$ws.Range("I22").Value = 1144.1052
$ws.Range("J22").Value = 4374.9165
$ws.Range("K22").Value = 1144.1052
$ws.Range("L22").Value = 4374.9165
$ws.Range("M22").Value = -849.1052
$ws.Range("N22").Value = -4964.9165

$ws.Range("H27").Value = 2394.742
$ws.Range("I27").Value = 1144.1052
$ws.Range("J27").Value = 4374.9165
$ws.Range("K27").Value = 1144.1052
$ws.Range("L27").Value = 4374.9165
$ws.Range("M27").Value = -1037.1052
$ws.Range("N27").Value = -4588.9165

$ws.Range("H40").Value = 10073.772
$ws.Range("I40").Value = 10057.637
$ws.Range("J40").Value = 10089.909
$ws.Range("K40").Value = 10057.637
$ws.Range("L40").Value = 10089.909
$ws.Range("M40").Value = -9921.637000000001
$ws.Range("N40").Value = -10361.909

$ws.Range("H48").Value = 10046
$ws.Range("J48").Value = 10046
$ws.Range("L48").Value = 10046
$ws.Range("N48").Value = -11368

$ws.Range("H93").Value = 1875.05
$ws.Range("I93").Value = 1842.3334
$ws.Range("K93").Value = 1842.3334
$ws.Range("M93").Value = -594.3334

$ws.Range("H126").Value = 10217.363
$ws.Range("I126").Value = 6848.6665
$ws.Range("J126").Value = 14259.8
$ws.Range("K126").Value = 20545.9995
$ws.Range("L126").Value = 42779.39999999999
$ws.Range("M126").Value = -18075.9995
$ws.Range("N126").Value = -47719.39999999999

$ws.Range("H132").Value = 4921.83
$ws.Range("I132").Value = 5013.6587
$ws.Range("K132").Value = 15040.9761
$ws.Range("M132").Value = -12510.9761

Write-Output "Applied all market-data updates."